# Auto-update draw results: append the 2025-09-26 Pick 3 draw as a new row.
#
# All values in this table are stored as plain text (t="str" in the OOXML),
# even though some of them look like dates/numbers (e.g. "2025-09-26",
# "250926"). To reproduce that faithfully we briefly force a text number
# format before assigning the value (so Excel doesn't silently convert the
# string into a date serial / numeric value), then restore the cell's style
# back to Normal so the new row doesn't end up with a lingering explicit
# cell style that wasn't part of the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 10

function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item($newRow, 1) "2025-09-26"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
Set-TextValue $ws.Cells.Item($newRow, 3) "250926"
$ws.Cells.Item($newRow, 4).Value = "1-1-0"
Set-TextValue $ws.Cells.Item($newRow, 5) "2025-09-26T21:37:22.557+04:00"
